$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns before the existing "budget"/"userId" columns, which
# shifts them from B,C to F,G.
$ws.Range("B1:E1").EntireColumn.Insert()

# The newly inserted cells under rows 2-3 inherited column A's bordered/bold
# style from the insert; these are plain data cells, so strip that back off.
$ws.Range("B2:E3").ClearFormats()

# Give the 4 new header cells the same header style as the existing headers.
$ws.Range("F1").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)

# New header row text
$ws.Range("B1").Value = "Unnamed: 0"
$ws.Range("C1").Value = "Unnamed: 0.1"
$ws.Range("D1").Value = "Unnamed: 0.1.1"
$ws.Range("E1").Value = "Unnamed: 0.1.1.1"

# Fill in the new index-like columns for the existing rows (2-3)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# New rows 4-7 (new_user window), with a "staircase" of trailing blanks
# (the blank B/C/D/E cells are left untouched - a never-written cell and a
# cell explicitly set to "" both serialise the same way, as empty/absent).
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = 3000

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 3
$ws.Range("F5").Value = 0

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 4
$ws.Range("F6").Value = 3000

$ws.Range("A7").Value = 5
$ws.Range("F7").Value = 0

# Random-looking userId values on the new rows must stay TEXT (not get
# auto-converted to numbers). Force text type via a text NumberFormat, then
# clear the format back off so no residual style is left on the cells.
$ws.Range("G4:G7").NumberFormat = "@"
$ws.Range("G4").Value = "0.18283207537334212"
$ws.Range("G5").Value = "0.9025862388675534"
$ws.Range("G6").Value = "0.1105236658581672"
$ws.Range("G7").Value = "0.09698904686853815"
$ws.Range("G4:G7").ClearFormats()

# Apply the bordered/bold/centered index style to the new index cells in
# column A (same style as A2:A3).
$ws.Range("A2").Copy()
$ws.Range("A4:A7").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
